$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits at the very
#    start of the document (before the "Name:" run in paragraph 1).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Merge the two runs "Lab Exercise 12.20" + ".2019" in paragraph 3
#    into a single run reading "Lab Exercise 12.23.2020".
# ---------------------------------------------------------------------
$findRange = $d.Range(0, $d.Content.End)
$findRange.Find.Execute("Lab Exercise 12.20.2019", $true, $false, $false, $false, $false, $true, 1, $false, "Lab Exercise 12.23.2020", 2)

# ---------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark, now collapsed right after the new
#    text, at the end of paragraph 3 (before its paragraph mark).
#
#    A bookmark range collapsed exactly on a paragraph mark can't be
#    created directly and reliably, so a one-character placeholder is
#    inserted at the target spot, the bookmark is anchored around it,
#    and then the placeholder is deleted - leaving the bookmark
#    collapsed at the correct position.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$endPos = $p3.Range.End - 1

$insertRange = $d.Range($endPos, $endPos)
$insertRange.InsertAfter("X")

$markerRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($endPos, $endPos + 1)
$markerRange2.Text = ""
